{"js": "// Apply the diff: replace each old math-fact / date string with its new value.\n// Every \"old\" text below is unique within the document body, so an exact,\n// case-sensitive whole-text search uniquely identifies the run to update.\nconst pairs = [[\"2024-03-16 Saturday\", \"2024-03-17 Sunday\"], [\"80-28=52\", \"67-56=11\"], [\"20+44=64\", \"15-8=7\"], [\"51-20=31\", \"95-69=26\"], [\"72-29=43\", \"76+21=97\"], [\"48+1=49\", \"80+5=85\"], [\"11+23=34\", \"1+53=54\"], [\"49-46=3\", \"41+22=63\"], [\"32-18=14\", \"40-26=14\"], [\"14+68=82\", \"58-15=43\"], [\"96-53=43\", \"4+0=4\"], [\"37-16=21\", \"20+33=53\"], [\"43+51=94\", \"44+31=75\"], [\"52+9=61\", \"3+13=16\"], [\"54-16=38\", \"16+7=23\"], [\"62-18=44\", \"98-66=32\"], [\"46+32=78\", \"70-37=33\"], [\"36+13=49\", \"50-8=42\"], [\"62+30=92\", \"70-61=9\"], [\"70-69=1\", \"29-12=17\"], [\"57-0=57\", \"24+39=63\"], [\"71+20=91\", \"61-23=38\"], [\"17+4=21\", \"2+70=72\"], [\"83-72=11\", \"69+20=89\"], [\"84-41=43\", \"7+47=54\"], [\"31+27=58\", \"45-27=18\"], [\"70-41=29\", \"64+29=93\"], [\"89+10=99\", \"41+12=53\"], [\"96-51=45\", \"34+47=81\"], [\"94-54=40\", \"83-54=29\"], [\"19+70=89\", \"21+30=51\"], [\"85-63=22\", \"56-9=47\"], [\"69+13=82\", \"50-43=7\"], [\"96-57=39\", \"90+4=94\"], [\"24+58=82\", \"51+36=87\"], [\"34+38=72\", \"44+52=96\"], [\"23+63=86\", \"72-10=62\"], [\"34+0=34\", \"36+34=70\"], [\"10+66=76\", \"96-78=18\"], [\"3+47=50\", \"44+0=44\"], [\"29+23=52\", \"98-25=73\"], [\"31+14=45\", \"77-18=59\"], [\"40-27=13\", \"8+31=39\"], [\"74+15=89\", \"4+2=6\"], [\"68+17=85\", \"62-56=6\"], [\"39+18=57\", \"26+9=35\"], [\"22+27=49\", \"43-35=8\"], [\"21+72=93\", \"46+12=58\"], [\"39+2=41\", \"75-70=5\"], [\"60-11=49\", \"84-76=8\"], [\"44+23=67\", \"18+6=24\"], [\"4+33=37\", \"99-23=76\"], [\"62+22=84\", \"76-14=62\"], [\"39-24=15\", \"70-65=5\"], [\"35+52=87\", \"98-44=54\"], [\"42+8=50\", \"84-64=20\"], [\"82+6=88\", \"35-3=32\"], [\"53+4=57\", \"58-52=6\"], [\"20+28=48\", \"8+48=56\"], [\"85-51=34\", \"13+41=54\"], [\"32+48=80\", \"28-6=22\"], [\"62+17=79\", \"2+65=67\"], [\"99-91=8\", \"45+6=51\"], [\"44+9=53\", \"50-2=48\"], [\"31+61=92\", \"97-4=93\"], [\"98-51=47\", \"75+15=90\"], [\"95-68=27\", \"55+9=64\"], [\"72-13=59\", \"16+44=60\"], [\"88-65=23\", \"82-74=8\"], [\"67-39=28\", \"17+14=31\"], [\"51+19=70\", \"21+60=81\"], [\"34+10=44\", \"30+18=48\"], [\"39+15=54\", \"18-14=4\"], [\"28+20=48\", \"62-58=4\"], [\"12+65=77\", \"34+16=50\"], [\"67+10=77\", \"66-18=48\"], [\"21+45=66\", \"11+55=66\"], [\"87+5=92\", \"43-23=20\"], [\"74-27=47\", \"48+28=76\"], [\"27-3=24\", \"52+8=60\"], [\"50-18=32\", \"50-22=28\"], [\"26+29=55\", \"64-46=18\"], [\"26+71=97\", \"10+4=14\"], [\"95-49=46\", \"32-26=6\"], [\"4+26=30\", \"26-9=17\"], [\"94-53=41\", \"29+64=93\"], [\"16+49=65\", \"25+70=95\"], [\"47+36=83\", \"31+30=61\"], [\"76-48=28\", \"72-17=55\"], [\"80-52=28\", \"9+52=61\"], [\"63-39=24\", \"31-29=2\"], [\"99-54=45\", \"78+5=83\"], [\"64-63=1\", \"90-68=22\"], [\"9+10=19\", \"79-21=58\"], [\"36+63=99\", \"58-0=58\"], [\"58-14=44\", \"62+32=94\"], [\"53+16=69\", \"38+23=61\"], [\"75-4=71\", \"99-72=27\"], [\"81+3=84\", \"44+20=64\"], [\"12+16=28\", \"61-25=36\"], [\"12+34=46\", \"23+65=88\"]];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  // Replace the first (and only) match's text in place so formatting (font,\n  // size, etc.) on the existing run is preserved.\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the diff: replace each old math-fact / date string with its new value.\n# Every \"old\" text below is unique within the document body, so an exact,\n# case-sensitive whole-text Find/Replace uniquely targets the run to update.\n$pairs = @(\n    @('2024-03-16 Saturday', '2024-03-17 Sunday'),\n    @('80-28=52', '67-56=11'),\n    @('20+44=64', '15-8=7'),\n    @('51-20=31', '95-69=26'),\n    @('72-29=43', '76+21=97'),\n    @('48+1=49', '80+5=85'),\n    @('11+23=34', '1+53=54'),\n    @('49-46=3', '41+22=63'),\n    @('32-18=14', '40-26=14'),\n    @('14+68=82', '58-15=43'),\n    @('96-53=43', '4+0=4'),\n    @('37-16=21', '20+33=53'),\n    @('43+51=94', '44+31=75'),\n    @('52+9=61', '3+13=16'),\n    @('54-16=38', '16+7=23'),\n    @('62-18=44', '98-66=32'),\n    @('46+32=78', '70-37=33'),\n    @('36+13=49', '50-8=42'),\n    @('62+30=92', '70-61=9'),\n    @('70-69=1', '29-12=17'),\n    @('57-0=57', '24+39=63'),\n    @('71+20=91', '61-23=38'),\n    @('17+4=21', '2+70=72'),\n    @('83-72=11', '69+20=89'),\n    @('84-41=43', '7+47=54'),\n    @('31+27=58', '45-27=18'),\n    @('70-41=29', '64+29=93'),\n    @('89+10=99', '41+12=53'),\n    @('96-51=45', '34+47=81'),\n    @('94-54=40', '83-54=29'),\n    @('19+70=89', '21+30=51'),\n    @('85-63=22', '56-9=47'),\n    @('69+13=82', '50-43=7'),\n    @('96-57=39', '90+4=94'),\n    @('24+58=82', '51+36=87'),\n    @('34+38=72', '44+52=96'),\n    @('23+63=86', '72-10=62'),\n    @('34+0=34', '36+34=70'),\n    @('10+66=76', '96-78=18'),\n    @('3+47=50', '44+0=44'),\n    @('29+23=52', '98-25=73'),\n    @('31+14=45', '77-18=59'),\n    @('40-27=13', '8+31=39'),\n    @('74+15=89', '4+2=6'),\n    @('68+17=85', '62-56=6'),\n    @('39+18=57', '26+9=35'),\n    @('22+27=49', '43-35=8'),\n    @('21+72=93', '46+12=58'),\n    @('39+2=41', '75-70=5'),\n    @('60-11=49', '84-76=8'),\n    @('44+23=67', '18+6=24'),\n    @('4+33=37', '99-23=76'),\n    @('62+22=84', '76-14=62'),\n    @('39-24=15', '70-65=5'),\n    @('35+52=87', '98-44=54'),\n    @('42+8=50', '84-64=20'),\n    @('82+6=88', '35-3=32'),\n    @('53+4=57', '58-52=6'),\n    @('20+28=48', '8+48=56'),\n    @('85-51=34', '13+41=54'),\n    @('32+48=80', '28-6=22'),\n    @('62+17=79', '2+65=67'),\n    @('99-91=8', '45+6=51'),\n    @('44+9=53', '50-2=48'),\n    @('31+61=92', '97-4=93'),\n    @('98-51=47', '75+15=90'),\n    @('95-68=27', '55+9=64'),\n    @('72-13=59', '16+44=60'),\n    @('88-65=23', '82-74=8'),\n    @('67-39=28', '17+14=31'),\n    @('51+19=70', '21+60=81'),\n    @('34+10=44', '30+18=48'),\n    @('39+15=54', '18-14=4'),\n    @('28+20=48', '62-58=4'),\n    @('12+65=77', '34+16=50'),\n    @('67+10=77', '66-18=48'),\n    @('21+45=66', '11+55=66'),\n    @('87+5=92', '43-23=20'),\n    @('74-27=47', '48+28=76'),\n    @('27-3=24', '52+8=60'),\n    @('50-18=32', '50-22=28'),\n    @('26+29=55', '64-46=18'),\n    @('26+71=97', '10+4=14'),\n    @('95-49=46', '32-26=6'),\n    @('4+26=30', '26-9=17'),\n    @('94-53=41', '29+64=93'),\n    @('16+49=65', '25+70=95'),\n    @('47+36=83', '31+30=61'),\n    @('76-48=28', '72-17=55'),\n    @('80-52=28', '9+52=61'),\n    @('63-39=24', '31-29=2'),\n    @('99-54=45', '78+5=83'),\n    @('64-63=1', '90-68=22'),\n    @('9+10=19', '79-21=58'),\n    @('36+63=99', '58-0=58'),\n    @('58-14=44', '62+32=94'),\n    @('53+16=69', '38+23=61'),\n    @('75-4=71', '99-72=27'),\n    @('81+3=84', '44+20=64'),\n    @('12+16=28', '61-25=36'),\n    @('12+34=46', '23+65=88')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop: do not wrap past the end of the content\n\n    # FindReplace signature:\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # wdReplace.wdReplaceAll = 2 (replace every match \u2014 each text is unique, so this\n    # replaces exactly the one run that should change).\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
